# Reproduce the manual edit: three stray/duplicate salary entries (rows 14,
# 26 and 35 of column B) were deleted from the sheet, and the window was
# left scrolled down with B35 (the last cell touched) selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Remove the three erroneous Salary values.
$ws.Range("B14").ClearContents()
$ws.Range("B26").ClearContents()
$ws.Range("B35").ClearContents()

# Leave the view scrolled so row 23 is at the top, with B35 as the active
# selection (the last cell edited).
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 23
$ws.Range("B35").Select()
